{"js": "// Replace the two-digit-number-divided-by-one-digit-number problems that\n// changed between before.docx and after.docx. Each \"old\" expression is\n// unique within the document, so a literal body.search() + replace is\n// unambiguous and safe.\nconst replacements = [\n  [\"67\u00f77=\", \"46\u00f76=\"],\n  [\"45\u00f73=\", \"19\u00f74=\"],\n  [\"51\u00f72=\", \"92\u00f76=\"],\n  [\"96\u00f78=\", \"74\u00f78=\"],\n  [\"17\u00f74=\", \"56\u00f77=\"],\n  [\"70\u00f79=\", \"53\u00f72=\"],\n  [\"35\u00f76=\", \"30\u00f78=\"],\n  [\"83\u00f73=\", \"96\u00f73=\"],\n  [\"56\u00f75=\", \"91\u00f74=\"],\n  [\"32\u00f77=\", \"88\u00f77=\"],\n  [\"32\u00f78=\", \"16\u00f73=\"],\n  [\"47\u00f78=\", \"42\u00f78=\"],\n  [\"62\u00f75=\", \"98\u00f74=\"],\n  [\"96\u00f72=\", \"59\u00f73=\"],\n  [\"90\u00f72=\", \"29\u00f78=\"],\n  [\"35\u00f72=\", \"66\u00f76=\"],\n  [\"83\u00f75=\", \"49\u00f73=\"],\n  [\"75\u00f75=\", \"32\u00f75=\"],\n  [\"64\u00f76=\", \"31\u00f72=\"],\n  [\"31\u00f73=\", \"60\u00f79=\"],\n  [\"57\u00f73=\", \"29\u00f73=\"],\n  [\"36\u00f77=\", \"45\u00f77=\"],\n  [\"49\u00f78=\", \"54\u00f79=\"],\n  [\"18\u00f75=\", \"94\u00f78=\"],\n  [\"60\u00f74=\", \"77\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the ten two-digit-number-divided-by-one-digit-number problems\n# that changed between before.docx and after.docx. Each \"old\" expression\n# is unique within the document, so Find/Replace (one hit per pair) is\n# unambiguous and safe.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"67\u00f77=\", \"46\u00f76=\"),\n    @(\"45\u00f73=\", \"19\u00f74=\"),\n    @(\"51\u00f72=\", \"92\u00f76=\"),\n    @(\"96\u00f78=\", \"74\u00f78=\"),\n    @(\"17\u00f74=\", \"56\u00f77=\"),\n    @(\"70\u00f79=\", \"53\u00f72=\"),\n    @(\"35\u00f76=\", \"30\u00f78=\"),\n    @(\"83\u00f73=\", \"96\u00f73=\"),\n    @(\"56\u00f75=\", \"91\u00f74=\"),\n    @(\"32\u00f77=\", \"88\u00f77=\"),\n    @(\"32\u00f78=\", \"16\u00f73=\"),\n    @(\"47\u00f78=\", \"42\u00f78=\"),\n    @(\"62\u00f75=\", \"98\u00f74=\"),\n    @(\"96\u00f72=\", \"59\u00f73=\"),\n    @(\"90\u00f72=\", \"29\u00f78=\"),\n    @(\"35\u00f72=\", \"66\u00f76=\"),\n    @(\"83\u00f75=\", \"49\u00f73=\"),\n    @(\"75\u00f75=\", \"32\u00f75=\"),\n    @(\"64\u00f76=\", \"31\u00f72=\"),\n    @(\"31\u00f73=\", \"60\u00f79=\"),\n    @(\"57\u00f73=\", \"29\u00f73=\"),\n    @(\"36\u00f77=\", \"45\u00f77=\"),\n    @(\"49\u00f78=\", \"54\u00f79=\"),\n    @(\"18\u00f75=\", \"94\u00f78=\"),\n    @(\"60\u00f74=\", \"77\u00f79=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    # wdReplaceOne = 1 -> replace just the single (unique) occurrence.\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n    if (-not $found) {\n        throw \"Could not find text '$oldText' to replace.\"\n    }\n}\n"}
